# Updates cryptos list (Price / Volume(1h) columns, and row 51 coin swap)
# applied via Excel COM interop, preserving text cell typing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "42.924.64"
Set-TextValue "E2" "  -0.06%  "
Set-TextValue "D3" "2.548.68"
Set-TextValue "E3" "  +0.19%  "
Set-TextValue "E4" "  -0.07%  "
Set-TextValue "D5" "303.99"
Set-TextValue "E5" "  +1.85%  "
Set-TextValue "D6" "98.00"
Set-TextValue "E6" "  +4.19%  "
Set-TextValue "D7" "0.578"
Set-TextValue "E7" "  +0.85%  "
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "D9" "0.545"
Set-TextValue "E9" "  -0.54%  "
Set-TextValue "D10" "36.66"
Set-TextValue "E10" "  +1.98%  "
Set-TextValue "D11" "0.0827"
Set-TextValue "E11" "  +3.08%  "
Set-TextValue "E12" "  +3.80%  "
Set-TextValue "E13" "  +0.77%  "
Set-TextValue "D14" "2.939.98"
Set-TextValue "E14" "  +0.03%  "
Set-TextValue "D15" "2.533.74"
Set-TextValue "E15" "  -0.68%  "
Set-TextValue "D16" "15.12"
Set-TextValue "E16" "  +7.53%  "
Set-TextValue "D17" "0.874"
Set-TextValue "E17" "  +0.79%  "
Set-TextValue "D18" "42.958.30"
Set-TextValue "E18" "  -0.18%  "
Set-TextValue "D19" "13.92"
Set-TextValue "E19" "  +7.27%  "
Set-TextValue "D20" "0.0₃0993"
Set-TextValue "E20" "  +1.33%  "
Set-TextValue "E21" "  -0.06%  "
Set-TextValue "D22" "71.96"
Set-TextValue "E22" "  +0.48%  "
Set-TextValue "D23" "254.52"
Set-TextValue "E23" "  -0.50%  "
Set-TextValue "E24" "  +1.98%  "
Set-TextValue "E25" "  -1.51%  "
Set-TextValue "D26" "27.99"
Set-TextValue "E26" "  -3.62%  "
Set-TextValue "E27" "  +0.02%  "
Set-TextValue "D28" "10.28"
Set-TextValue "E28" "  +2.83%  "
Set-TextValue "D29" "37.72"
Set-TextValue "E29" "  +1.69%  "
Set-TextValue "E30" "  -1.65%  "
Set-TextValue "D31" "6.16"
Set-TextValue "E31" "  +4.18%  "
Set-TextValue "D32" "158.73"
Set-TextValue "E32" "  +4.06%  "
Set-TextValue "D33" "19.40"
Set-TextValue "E33" "  +14.57%  "
Set-TextValue "D34" "2.14"
Set-TextValue "E34" "  -0.51%  "
Set-TextValue "E35" "  +1.17%  "
Set-TextValue "E36" "  -1.96%  "
Set-TextValue "E37" "  -4.48%  "
Set-TextValue "E38" "  +2.32%  "
Set-TextValue "D39" "25.42"
Set-TextValue "E39" "  +9.68%  "
Set-TextValue "E40" "  -0.07%  "
Set-TextValue "D41" "2.10"
Set-TextValue "E41" "  +31.90%  "
Set-TextValue "E42" "  +0.04%  "
Set-TextValue "D43" "3.90"
Set-TextValue "E43" "  +0.32%  "
Set-TextValue "D44" "2.096.37"
Set-TextValue "E44" "  +0.63%  "
Set-TextValue "D45" "0.0307"
Set-TextValue "E45" "  -1.22%  "
Set-TextValue "D46" "0.998"
Set-TextValue "E46" "  -0.13%  "
Set-TextValue "D47" "86.37"
Set-TextValue "E47" "  +2.52%  "
Set-TextValue "D48" "8.93"
Set-TextValue "E48" "  +0.67%  "
Set-TextValue "D49" "75.41"
Set-TextValue "E49" "  +9.93%  "
Set-TextValue "D50" "2.798.28"
Set-TextValue "E50" "  +0.05%  "
Set-TextValue "B51" "Algorand"
Set-TextValue "C51" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D51" "0.192"
Set-TextValue "E51" "  +3.03%  "
